$wb = $excel.ActiveWorkbook

# --- 1. METADATA_MEASURES: update description of MQME018 (cell C8) ---
$wsMeasures = $wb.Worksheets.Item("METADATA_MEASURES")
$wsMeasures.Range("C8").Value = "Total number of cells in schema (sum of columns x rows for each table)"

# --- 2. METADATA_ISSUES: add SUGGESTED_* columns (I:L) with suggestion data ---
$wsIssues = $wb.Worksheets.Item("METADATA_ISSUES")

# Header row
$wsIssues.Range("I1").Value = "SUGGESTED_VALUE"
$wsIssues.Range("J1").Value = "SUGGESTED_SOURCE"
$wsIssues.Range("K1").Value = "SUGGESTED_CONFIDENCE"
$wsIssues.Range("L1").Value = "SUGGESTED_DDL"
# Copy the header formatting (bold, border, centered) from an existing header cell
$wsIssues.Range("A1").Copy()
$wsIssues.Range("I1:L1").PasteSpecial(-4122)

# Row 2 - PEDIDO.UID_GERACAO_DEBITO (non-standard prefix) -> suggest rename
$wsIssues.Range("I2").Value = "DSC_GERACAO_DEBITO"
$wsIssues.Range("J2").Value = "RULES"
$wsIssues.Range("K2").Value = 0.85
$wsIssues.Range("L2").Value = "ALTER TABLE SISAGUA.PEDIDO RENAME COLUMN UID_GERACAO_DEBITO TO DSC_GERACAO_DEBITO;"

# Row 3 - GRAFICA.SEQ_GRAFICA -> suggest comment
$wsIssues.Range("I3").Value = "Sequencia de grafica."
$wsIssues.Range("J3").Value = "RULES"
$wsIssues.Range("K3").Value = 0.8
$wsIssues.Range("L3").Value = "COMMENT ON COLUMN SISAGUA.GRAFICA.SEQ_GRAFICA IS 'Sequencia de grafica.';"

# Row 4 - GRAFICA.COD_CNPJ -> suggest comment
$wsIssues.Range("I4").Value = "Codigo de cnpj."
$wsIssues.Range("J4").Value = "RULES"
$wsIssues.Range("K4").Value = 0.8
$wsIssues.Range("L4").Value = "COMMENT ON COLUMN SISAGUA.GRAFICA.COD_CNPJ IS 'Codigo de cnpj.';"

# Row 5 - GRAFICA.NOM_GRAFICA -> suggest comment
$wsIssues.Range("I5").Value = "Nome de grafica."
$wsIssues.Range("J5").Value = "RULES"
$wsIssues.Range("K5").Value = 0.8
$wsIssues.Range("L5").Value = "COMMENT ON COLUMN SISAGUA.GRAFICA.NOM_GRAFICA IS 'Nome de grafica.';"

# Row 6 - GRAFICA.DAT_INICIO -> suggest comment
$wsIssues.Range("I6").Value = "Data de inicio."
$wsIssues.Range("J6").Value = "RULES"
$wsIssues.Range("K6").Value = 0.8
$wsIssues.Range("L6").Value = "COMMENT ON COLUMN SISAGUA.GRAFICA.DAT_INICIO IS 'Data de inicio.';"

# Row 7 - GRAFICA.DAT_FIM -> suggest comment
$wsIssues.Range("I7").Value = "Data de fim."
$wsIssues.Range("J7").Value = "RULES"
$wsIssues.Range("K7").Value = 0.8
$wsIssues.Range("L7").Value = "COMMENT ON COLUMN SISAGUA.GRAFICA.DAT_FIM IS 'Data de fim.';"

# Row 8 - GRAFICA.STA_ATIVA -> suggest comment
$wsIssues.Range("I8").Value = "Status de ativa."
$wsIssues.Range("J8").Value = "RULES"
$wsIssues.Range("K8").Value = 0.8
$wsIssues.Range("L8").Value = "COMMENT ON COLUMN SISAGUA.GRAFICA.STA_ATIVA IS 'Status de ativa.';"

# Row 9 - GRAFICA.NOM_USUARIO -> suggest comment
$wsIssues.Range("I9").Value = "Nome de usuario."
$wsIssues.Range("J9").Value = "RULES"
$wsIssues.Range("K9").Value = 0.8
$wsIssues.Range("L9").Value = "COMMENT ON COLUMN SISAGUA.GRAFICA.NOM_USUARIO IS 'Nome de usuario.';"

# Row 10 - GRAFICA.NOM_SENHA -> suggest comment
$wsIssues.Range("I10").Value = "Nome de senha."
$wsIssues.Range("J10").Value = "RULES"
$wsIssues.Range("K10").Value = 0.8
$wsIssues.Range("L10").Value = "COMMENT ON COLUMN SISAGUA.GRAFICA.NOM_SENHA IS 'Nome de senha.';"

# Row 11 - GRAFICA.NOM_URL -> suggest comment
$wsIssues.Range("I11").Value = "Nome de url."
$wsIssues.Range("J11").Value = "RULES"
$wsIssues.Range("K11").Value = 0.8
$wsIssues.Range("L11").Value = "COMMENT ON COLUMN SISAGUA.GRAFICA.NOM_URL IS 'Nome de url.';"

# Row 12 - PEDIDO.SEQ_GRAFICA -> suggest comment
$wsIssues.Range("I12").Value = "Sequencia de grafica."
$wsIssues.Range("J12").Value = "RULES"
$wsIssues.Range("K12").Value = 0.8
$wsIssues.Range("L12").Value = "COMMENT ON COLUMN SISAGUA.PEDIDO.SEQ_GRAFICA IS 'Sequencia de grafica.';"

# Row 13 - PEDIDO_GRAFICA.SEQ_GRAFICA -> suggest comment
$wsIssues.Range("I13").Value = "Sequencia de grafica."
$wsIssues.Range("J13").Value = "RULES"
$wsIssues.Range("K13").Value = 0.8
$wsIssues.Range("L13").Value = "COMMENT ON COLUMN SISAGUA.PEDIDO_GRAFICA.SEQ_GRAFICA IS 'Sequencia de grafica.';"
